$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RuntimesChart")

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 0.0018369199999999999

$ws.Range("A8:B8").Select()
